$wb = $excel.ActiveWorkbook

$motorCalcs = $wb.Worksheets.Item("motorCalcs")
$motorCalcs.Range("B2").Value = 0.93720000000000003

$batteryCalcs = $wb.Worksheets.Item("batteryCalcs")
$batteryCalcs.Range("B2").Value = 11.990399999999999
